# Updates the crypto price/volume table (columns D and E) with refreshed
# values, matching the "Updated cryptos list ... with GitHub Actions" commit.
#
# Column D cells whose new value is a plain decimal number are first marked
# as Text (NumberFormat "@") before the assignment so Excel keeps them as
# strings (matching the source data's inlineStr cell type) instead of
# silently parsing them into numeric values; the style is then reset back
# to Normal so no residual cell formatting is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.242.02"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "3.031.06"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.18%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.029.01"
$ws.Range("E8").Value = "  +1.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.19%  "
$ws.Range("E11").Value = "  -1.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.489"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000248"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.30%  "
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").Value = "66.308.93"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").Value = "3.532.99"
$ws.Range("E17").Value = "  +1.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +19.39%  "
$ws.Range("D20").Value = "3.030.83"
$ws.Range("E20").Value = "  +1.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "466.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.09%  "
$ws.Range("E22").Value = "  +3.66%  "
$ws.Range("E23").Value = "  +0.79%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.61%  "
$ws.Range("E26").Value = "  -1.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.39%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  +1.05%  "
$ws.Range("E30").Value = "  +0.95%  "
$ws.Range("E31").Value = "  +1.20%  "
$ws.Range("E32").Value = "  +6.72%  "
$ws.Range("E33").Value = "  -3.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.85"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.989"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.05"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "49.51"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E42").Value = "  -1.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.65"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0360"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "378.87"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.87%  "
$ws.Range("D47").Value = "2.703.14"
$ws.Range("E47").Value = "  -2.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.55%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.43"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.37%  "
$ws.Range("E51").Value = "  +4.15%  "
